# Add MX-BBX and MX-DPBX accessories to the panel accessories sheets for
# Austria, Italy, Slovakia, Netherlands and Denmark markets.
#
# In each of those sheets the accessory list ends with two fixed rows
# ("Wg" then "Accessories"). The new rows are inserted right above the
# "Wg" row, pushing it (and the trailing "Accessories" row) down by two
# rows, while the two brand-new rows get the plain formatting used by
# the rest of the accessory list (copied from the row right above the
# insertion point).

$wb = $excel.ActiveWorkbook

$targets = @("Slovakia", "Italy", "Netherlands", "Austria", "Denmark")

foreach ($name in $targets) {
    $ws = $wb.Worksheets.Item($name)

    # Locate the last row of the sheet ("Accessories") - the row right
    # above it is the "Wg" row; new rows get inserted right before it.
    $lastRow = $ws.UsedRange.Rows.Count
    $insertRow = $lastRow - 1

    # Insert two blank rows, shifting "Wg"/"Accessories" rows down.
    $rowRange = $ws.Range($ws.Cells.Item($insertRow, 1), $ws.Cells.Item($insertRow + 1, 1)).EntireRow
    $rowRange.Insert()

    # Populate the new rows with the accessory names.
    $ws.Cells.Item($insertRow, 1).Value = "MX-DPBX"
    $ws.Cells.Item($insertRow + 1, 1).Value = "MX-BBX"

    # Copy formatting (style/border) from the row just above so the new
    # cells match the rest of the accessory list.
    $ws.Cells.Item($insertRow - 1, 1).Copy()
    $ws.Cells.Item($insertRow, 1).PasteSpecial(-4122)
    $ws.Cells.Item($insertRow + 1, 1).PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    # Update the selection to highlight the newly added rows.
    $ws.Activate()
    [void]$ws.Range($ws.Cells.Item($insertRow, 1), $ws.Cells.Item($insertRow + 1, 1)).Select()
}

# Make Slovakia the active tab, matching the saved workbook view state.
$wb.Worksheets.Item("Slovakia").Activate()
